$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 365
    3  = 366
    4  = 368
    5  = 370
    6  = 372
    7  = 374
    8  = 375
    9  = 377
    10 = 379
    11 = 381
    12 = 383
    13 = 386
    14 = 388
    15 = 33
    16 = 52
    17 = 103
    18 = 148
    19 = 172
    20 = 228
    21 = 248
    22 = 272
    23 = 315
    24 = 413
    25 = 421
    26 = 477
    27 = 514
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
